# Split the combined DoseTime (date+time) column into a date-only
# column (A) and a new time-only column (D, header "time").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for the time-of-day column.
$ws.Range("D1").Value = "time"

# Date-only serials (was date+time combined) for column A.
$dates = @(45623, 45624, 45625, 45626, 45627, 45628, 45629, 45630, 45631)

# Time-of-day fractions pulled out of the old column A values, now in D.
$times = @(
    0.375,
    0.54166666666666663,
    0.75,
    0.375,
    0.54166666666666663,
    0.75,
    0.375,
    0.54166666666666663,
    0.75
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 4).Value = $times[$i]
}

# Apply the date-only display format to column A, using a single
# source cell + copy/paste-format so the whole range shares one style
# instead of getting a distinct style per cell.
$ws.Range("A2").NumberFormat = "mm-dd-yy"
$ws.Range("A2").Copy()
$ws.Range("A3:A10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Apply the time-only display format to the new column D the same way.
$ws.Range("D2").NumberFormat = "[`$-F400]h:mm:ss AM/PM"
$ws.Range("D2").Copy()
$ws.Range("D3:D10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column B (PrescriptionScheduleEntry id) keeps its integer formatting.
$ws.Range("B2:B10").NumberFormat = "0"

$ws.Range("N14").Select()
